# Update crypto price/volume snapshot data (GitHub Actions scheduled refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.041.20"
$ws.Range("E2").Value = "  -1.49%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.666.21"
$ws.Range("E3").Value = "  -1.21%  "
$ws.Range("E4").Value = "  -0.13%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "215.99"
$ws.Range("E5").Value = "  -1.51%  "
$ws.Range("E6").Value = "  +0.67%  "
$ws.Range("E7").Value = "  -0.13%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2666"
$ws.Range("E8").Value = "  -0.12%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06390"
$ws.Range("E9").Value = "  +1.58%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "21.83"
$ws.Range("E10").Value = "  -0.88%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07453"
$ws.Range("E11").Value = "  +1.19%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.676.59"
$ws.Range("E12").Value = "  -0.60%  "
$ws.Range("E13").Value = "  -0.58%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.5803"
$ws.Range("E14").Value = "  +0.47%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.000008503"
$ws.Range("E15").Value = "  -0.88%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.04"
$ws.Range("E16").Value = "  -1.70%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "26.001.11"
$ws.Range("E17").Value = "  -1.87%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "4.922"
$ws.Range("E18").Value = "  -1.60%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.004"
$ws.Range("E19").Value = "  -0.18%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "10.78"
$ws.Range("E20").Value = "  -0.94%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "189.88"
$ws.Range("E21").Value = "  +2.31%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.184"
$ws.Range("E22").Value = "  -1.17%  "
$ws.Range("E23").Value = "  -0.07%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "144.89"
$ws.Range("E24").Value = "  +0.17%  "
$ws.Range("E25").Value = "  +1.56%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1212"
$ws.Range("E26").Value = "  +4.09%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "15.64"
$ws.Range("E27").Value = "  -0.23%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.06638"
$ws.Range("E28").Value = "  +15.67%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.328"
$ws.Range("E29").Value = "  -1.01%  "
$ws.Range("E30").Value = "  -1.79%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.548"
$ws.Range("E31").Value = "  +0.67%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.510"
$ws.Range("E32").Value = "  -0.24%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.660"
$ws.Range("E33").Value = "  +0.14%  "
$ws.Range("E34").Value = "  +0.03%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.6140"
$ws.Range("E35").Value = "  +3.20%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.370"
$ws.Range("E36").Value = "  +0.44%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.688"
$ws.Range("E37").Value = "  +0.32%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.380"
$ws.Range("E38").Value = "  +8.14%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.093.46"
$ws.Range("E39").Value = "  -0.78%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01591"
$ws.Range("E40").Value = "  -1.21%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.8695"
$ws.Range("E41").Value = "  +0.58%  "
$ws.Range("E42").Value = "  +0.49%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "101.53"
$ws.Range("E43").Value = "  +1.65%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.813.94"
$ws.Range("E44").Value = "  -1.68%  "
$ws.Range("E45").Value = "  -2.15%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "56.26"
$ws.Range("E46").Value = "  -0.08%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.008"
$ws.Range("E47").Value = "  +0.40%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.104"
$ws.Range("E48").Value = "  +0.80%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.05231"
$ws.Range("E49").Value = "  +0.35%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.4285"
$ws.Range("E50").Value = "  -0.68%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "5.990"
$ws.Range("E51").Value = "  +3.08%  "
